$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = '@'
$ws.Range("D2").Value = '45.419.94'
$ws.Range("D2").Style = 'Normal'
$ws.Range("E2").Value = '  -0.03%  '
$ws.Range("D3").NumberFormat = '@'
$ws.Range("D3").Value = '2.377.98'
$ws.Range("D3").Style = 'Normal'
$ws.Range("E3").Value = '  -0.29%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '316.42'
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  -0.64%  '
$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '108.77'
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  -3.49%  '
$ws.Range("D7").NumberFormat = '@'
$ws.Range("D7").Value = '0.641'
$ws.Range("D7").Style = 'Normal'
$ws.Range("E7").Value = '  +0.63%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -1.58%  '
$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '40.97'
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  -3.98%  '
$ws.Range("E11").Value = '  -1.53%  '
$ws.Range("E12").Value = '  -1.38%  '
$ws.Range("E13").Value = '  +0.88%  '
$ws.Range("E14").Value = '  -3.08%  '
$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '2.737.24'
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  -0.51%  '
$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '15.52'
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = '  -2.32%  '
$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '2.368.36'
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  -0.91%  '
$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '45.401.72'
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = '  +0.09%  '
$ws.Range("E19").Value = '  +21.96%  '
$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '7.36'
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  -3.73%  '
$ws.Range("E21").Value = '  -1.17%  '
$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '3.66'
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  +2.59%  '
$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '73.32'
$ws.Range("D23").Style = 'Normal'
$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '261.39'
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  -3.13%  '
$ws.Range("E25").Value = '  -0.98%  '
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '7.63'
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  +1.89%  '
$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '11.21'
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  -0.61%  '
$ws.Range("E29").Value = '  -1.57%  '
$ws.Range("E30").Value = '  +2.46%  '
$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '22.41'
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  -2.32%  '
$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '37.31'
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  -5.57%  '
$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '167.07'
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = '  -1.84%  '
$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '2.86'
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  -4.93%  '
$ws.Range("E35").Value = '  -0.04%  '
$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '0.117'
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  -0.40%  '
$ws.Range("E37").Value = '  -3.76%  '
$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '4.09'
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  +2.46%  '
$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '1.92'
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  +10.15%  '
$ws.Range("E40").Value = '  -2.16%  '
$ws.Range("E41").Value = '  -4.69%  '
$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '98.20'
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  -7.20%  '
$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '70.59'
$ws.Range("D43").Style = 'Normal'
$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '13.18'
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  -2.34%  '
$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '0.230'
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  -4.69%  '
$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '6.05'
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  +4.44%  '
$ws.Range("E47").Value = '  +0.08%  '
$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '1.832.91'
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  +11.34%  '
$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '84.32'
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  +7.52%  '
$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '112.29'
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  -4.36%  '
$ws.Range("E51").Value = '  -0.75%  '
